# Append 11 new homeworking-eligibility records (rows 87-97) to the
# "Original" sheet, matching the run made "for 0.85 sensitivity".
#
# Columns: A=primerNombre, B=segundoNombre, C=apellidoPaterno, D=apellidoMaterno

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("Luis",      "Alberto",    "Sotelo",    "Cordero"),
  @("Edwin",     "Ivan",       "Fernandez", "Castillo"),
  @("Jazmin",    "",           "Regalado",  "Cabello"),
  @("Percy",     "",           "Aguilar",   "Marin"),
  @("Christian", "Gianfranco", "Avalos",    "Cristobal"),
  @("Javier",    "Guillermo",  "Gonzales",  "Sandoval"),
  @("Diego",     "Joaquin",    "Guillen",   "Huarcaya"),
  @("Julio",     "",           "Arguedas",  "Quiñonez"),
  @("Danmert",   "Jonelly",    "Costilla",  "Claros"),
  @("Rodail",    "",           "Chavez",    "Rugel"),
  @("Kevin",     "Martin",     "Reyes",     "Reyes")
)

$startRow = 87
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $row = $data[$i]

  # The last populated row (86) carries the bordered / centered look used
  # throughout the table. Copy its formatting down into the new row first
  # (this reuses the existing cell styles instead of fabricating new ones),
  # then overwrite the values.
  $ws.Range("A86:D86").Copy()
  $ws.Range("A$r`:D$r").PasteSpecial(-4122)

  $ws.Cells.Item($r, 1).Value = $row[0]
  if ($row[1] -ne "") {
    $ws.Cells.Item($r, 2).Value = $row[1]
  }
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = $false

# Match the author's final selection/scroll position in the sheet.
$ws.Range("C53").Select() | Out-Null
